$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.ClearFormats()
}

$ws.Range("D2").Value = "35.574.93"
$ws.Range("E2").Value = "  +3.22%  "

$ws.Range("D3").Value = "1.841.75"
$ws.Range("E3").Value = "  +2.15%  "

$ws.Range("E4").Value = "  +0.19%  "

Set-TextValue "D5" "231.97"
$ws.Range("E5").Value = "  +3.33%  "

Set-TextValue "D6" "0.619"
$ws.Range("E6").Value = "  +2.79%  "

$ws.Range("E7").Value = "  +0.15%  "

Set-TextValue "D8" "43.74"
$ws.Range("E8").Value = "  +10.76%  "

Set-TextValue "D9" "0.311"
$ws.Range("E9").Value = "  +8.24%  "

Set-TextValue "D10" "0.0703"
$ws.Range("E10").Value = "  +5.21%  "

$ws.Range("E11").Value = "  +2.39%  "

$ws.Range("D12").Value = "2.109.35"
$ws.Range("E12").Value = "  +2.14%  "

$ws.Range("D13").Value = "1.849.17"
$ws.Range("E13").Value = "  +2.55%  "

Set-TextValue "D14" "11.32"
$ws.Range("E14").Value = "  +3.12%  "

Set-TextValue "D15" "0.675"
$ws.Range("E15").Value = "  +6.94%  "

Set-TextValue "D16" "4.74"
$ws.Range("E16").Value = "  +8.52%  "

$ws.Range("D17").Value = "35.548.11"
$ws.Range("E17").Value = "  +3.19%  "

Set-TextValue "D18" "70.42"
$ws.Range("E18").Value = "  +3.46%  "

$ws.Range("D19").Value = "0.0₃0802"
$ws.Range("E19").Value = "  +4.57%  "

Set-TextValue "D20" "245.07"
$ws.Range("E20").Value = "  +2.46%  "

Set-TextValue "D21" "12.04"
$ws.Range("E21").Value = "  +8.17%  "

Set-TextValue "D22" "4.64"
$ws.Range("E22").Value = "  +13.85%  "

$ws.Range("E23").Value = "  +0.19%  "

$ws.Range("E24").Value = "  +2.65%  "

Set-TextValue "D25" "172.29"
$ws.Range("E25").Value = "  +0.53%  "

Set-TextValue "D26" "7.98"
$ws.Range("E26").Value = "  +4.10%  "

Set-TextValue "D27" "17.85"
$ws.Range("E27").Value = "  +1.49%  "

Set-TextValue "D28" "0.122"
$ws.Range("E28").Value = "  +0.74%  "

Set-TextValue "D29" "1.55"
$ws.Range("E29").Value = "  +26.82%  "

$ws.Range("E30").Value = "  +0.18%  "

$ws.Range("D31").Value = "3.334.40"
$ws.Range("E31").Value = "  +37.24%  "

Set-TextValue "D32" "0.0552"
$ws.Range("E32").Value = "  +7.76%  "

$ws.Range("B33").Value = "InternetComputer(DFINITY)"
$ws.Range("C33").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue "D33" "4.09"
$ws.Range("E33").Value = "  +6.98%  "

$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue "D34" "3.94"
$ws.Range("E34").Value = "  +5.31%  "

Set-TextValue "D35" "1.84"

$ws.Range("E36").Value = "  +17.48%  "

Set-TextValue "D37" "1.14"
$ws.Range("E37").Value = "  +8.37%  "

Set-TextValue "D38" "0.690"
$ws.Range("E38").Value = "  +7.74%  "

$ws.Range("D39").Value = "1.349.04"
$ws.Range("E39").Value = "  +3.58%  "

$ws.Range("E40").Value = "  +7.23%  "

Set-TextValue "D41" "15.48"
$ws.Range("E41").Value = "  +11.02%  "

$ws.Range("E42").Value = "  +5.10%  "

$ws.Range("E43").Value = "  +7.35%  "

$ws.Range("E44").Value = "  +4.34%  "

$ws.Range("E45").Value = "  +0.65%  "

$ws.Range("E46").Value = "  +0.99%  "

$ws.Range("E47").Value = "  +9.96%  "

$ws.Range("E48").Value = "  +0.50%  "

$ws.Range("D49").Value = "2.013.38"
$ws.Range("E49").Value = "  +2.43%  "

$ws.Range("E50").Value = "  +0.21%  "

Set-TextValue "D51" "102.66"
$ws.Range("E51").Value = "  +0.77%  "
